# Extend the "Посыл / Заповедь" schedule table on Sheet1:
#  - column B (rows 2-11) all get the same (new) encrypted message text
#  - column C (rows 2-11) gets the updated list of time ranges, replacing
#    the old "18:55-19:0"/"19:0-19:5" rows with four new slots, and adding
#    two more rows at the end (21:25-21:30 / 21:30-21:35)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$posyl = 'U2FsdGVkX18TZnPw9mSh8RvSW4+879mztM+lFAIjdQNi5gBI8Te3ngM8KQbwAlXcm3N/M7mKtkzGzavZieqPMDoCvWWTxJ+O6t4ztmfvS7EmHZDLHw+wp6bspq2Yve7xnTxhx0sCw3nc4IWmsLvX3g=='

$times = @(
    '2:55-3:0',
    '3:0-3:5',
    '10:55-11:0',
    '11:0-11:5',
    '12:25-12:30',
    '12:50-12:55',
    '14:10-14:15',
    '14:15-14:20',
    '21:25-21:30',
    '21:30-21:35'
)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $times[$i]
}

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 2).Value = $posyl
}

for ($i = 8; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $times[$i]
}

$ws.Range("B16").Select() | Out-Null
